# Pythonlearn-02-Expressions.pptx
#
# The "Hello world" code sample on slide 2 reads:
#     >>> print('Hello world' )
# Fix the stray space before the closing paren so it reads:
#     >>> print('Hello world')
#
# We search every slide/shape for the text frame that contains the
# distinctive snippet rather than hard-coding slide/shape numbers, so the
# script is resilient to any incidental reshuffling.

$p = $ppt.ActivePresentation

$needle = "'Hello world' )"
$targetShape = $null
$targetRange = $null
$matchIndex0 = -1

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $range = $shape.TextFrame.TextRange
            $full = $range.Text
            if ($full.Contains($needle)) {
                $targetShape = $shape
                $targetRange = $range
                $matchIndex0 = $full.IndexOf($needle)
            }
        }
    }
}

if ($matchIndex0 -ge 0) {
    # 0-based offset of the ' )' (space + close-paren) inside the full text.
    $spaceParenOffset0 = $matchIndex0 + $needle.Length - 2
    # PowerPoint TextRange.Characters is 1-based.
    $start1 = $spaceParenOffset0 + 1

    $badBit = $targetRange.Characters($start1, 2)
    $badBit.Text = ")"
}
